# Insert a new weekly record at row 108 (Hortaliza, Vega Monumental Concepción - Berenjena).
# All existing records from row 108 down to row 141 shift down by one row
# (to rows 109-142) and the new record's data is written into row 108.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 108:141 down to 109:142, leaving row 108 empty for the new entry.
$ws.Rows.Item(108).Insert()

$ws.Range("A108").Value = 11
$ws.Range("B108").Value = "Vega Monumental Concepción"
$ws.Range("C108").Value = "Bíobío"
$ws.Range("D108").Value = 44988
$ws.Range("E108").Value = 8
$ws.Range("F108").Value = 100112001
$ws.Range("G108").Value = "Berenjena"
$ws.Range("H108").Value = "Sin especificar"
$ws.Range("I108").Value = "Primera"
$ws.Range("J108").Value = 270
$ws.Range("K108").Value = 7000
$ws.Range("L108").Value = 7500
$ws.Range("M108").Value = 7278
$ws.Range("N108").Value = "$/caja 60 unidades"
$ws.Range("O108").Value = "Región de Arica y Parinacota"
$ws.Range("P108").Value = 121
$ws.Range("Q108").Value = 60
$ws.Range("R108").Value = "Hortaliza"
